$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy H1's formatting, then set text)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2 through 36: column I is always 1, column J mirrors column H
for ($row = 2; $row -le 36; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValue
}
